$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# --- Update existing entries ---

# 11/3 (row 20): ARVO abstract talks -> now starts at 10am
$ws.Range("G20").Value = "ARVO abstract talks. Start at 10am"

# 11/17 (row 22): add lab meeting presenter (Justin) and food (Solomon)
$ws.Range("B22").Value = "Justin"
$ws.Range("E22").Value = "Solomon"

# --- Insert a new week (11/20) right after the 11/17 row ---
$ws.Rows("23:23").Insert()

# New row 23: 11/20/2025 - Solomon's PhD
$ws.Range("A23").Value = 45981
$ws.Range("F23").Formula = "=FALSE"
$ws.Range("B23").Value = "Solomon"
$ws.Range("G23").Value = "Solomon's PhD"

# --- Fill in notes for the following (now shifted) weeks ---

# 11/24 (row 24): no meeting - Thanksgiving
$ws.Range("G24").Value = "no meeting - Thanksgiving"

# 12/1 (row 25): Everyone - ARVO abstract editing
$ws.Range("B25").Value = "Everyone"
$ws.Range("G25").Value = "ARVO abstract editing"

# Leave the cursor where the author left it when saving
$ws.Range("C21").Select() | Out-Null
